$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 6 (pushes existing rows 6..84 down to 7..85,
# and extends the used range to A1:T85). The new row's header/layout comes
# from the row formatting already present (e.g. the date column's style).
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new weekly price record.
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C6").Value = "Arica y Parinacota"
$ws.Range("D6").Value = 44685
$ws.Range("E6").Value = 15
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100102
$ws.Range("H6").Value = "Cítricos"
$ws.Range("I6").Value = 100102005
$ws.Range("J6").Value = "Naranja"
$ws.Range("K6").Value = "Valencia"
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 300
$ws.Range("N6").Value = 950
$ws.Range("O6").Value = 1000
$ws.Range("P6").Value = 975
$ws.Range("Q6").Value = "`$/kilo (en caja de 20 kilos)"
$ws.Range("R6").Value = "Región de O'Higgins"
$ws.Range("S6").Value = 975
$ws.Range("T6").Value = 1
